# Update the "last updated" timestamp text in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 23:34"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1342430
$ws.Range("C4").Value = 20645
$ws.Range("D4").Value = 232827
$ws.Range("E4").Value = 1029686
$ws.Range("F4").Value = 16801
$ws.Range("G4").Value = 1302
$ws.Range("H4").Value = 79917

# Row 10 - Alemania
$ws.Range("B10").Value = 171324
$ws.Range("C10").Value = 736
$ws.Range("E10").Value = 20475
$ws.Range("G10").Value = 39
$ws.Range("H10").Value = 7549

# Row 11 - Brasil
$ws.Range("B11").Value = 149101
$ws.Range("C11").Value = 3209
$ws.Range("E11").Value = 79691
$ws.Range("G11").Value = 121
$ws.Range("H11").Value = 10113

# Row 78 - Guinea
$ws.Range("B78").Value = 2042
$ws.Range("C78").Value = 33
$ws.Range("D78").Value = 698
$ws.Range("E78").Value = 1333

# Row 140 - Cabo Verde
$ws.Range("D140").Value = 56
$ws.Range("E140").Value = 178

# Row 151 - Haiti
$ws.Range("B151").Value = 151
$ws.Range("C151").Value = 5
$ws.Range("E151").Value = 122
